$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Number" header column (I1), matching the bold/centered header style
# but using a text number format (so phone numbers aren't mangled) and no border.
$ws.Range("I1").Value = "Number"
$ws.Range("I1").Font.Bold = $true
$ws.Range("I1").HorizontalAlignment = -4108
$ws.Range("I1").VerticalAlignment = -4160

# Patient phone numbers
$ws.Range("I2").Value = 81234567
$ws.Range("I3").Value = 81112222
$ws.Range("I4").Value = 91234567

# Store the whole new column as text-formatted numbers
$ws.Columns("I").NumberFormat = "@"

# Column width adjustments
$ws.Columns("E").ColumnWidth = 16.666666666666668
$ws.Columns("H").ColumnWidth = 18
$ws.Columns("I").ColumnWidth = 8.25

# Update the active selection
$ws.Range("H7").Select() | Out-Null
